$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = "Evans Honkapohja 2001"
$ws.Range("B10").Value = "p. 41"
$ws.Range("C10").Value = "overparameterization of PLM"

$ws.Range("C11").Select()
